$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$articleText = @"
h1: Why are we divided?
p.note: We have completed Surah Baqarah today. in’sha’Allah with this pace we will complete Quran by December 2020.
p: Seriously, is there really a need to call myself a Sunni Muslim. Or Shia, Wahabi or Marzai. We have disintegrated ourselves into sects and created a divide within us, ourselves. If Europe can unite with all its religions and secularism why can’t we unite. 
p: Shia, Sunni, Wahabi, Dyobandi, Marzai are all names of groups with deviated intentions. A muslim, true muslim who lives his life by the book does not really care to which sect he belongs. He just calls himself a muslim. He lives by the principles written in one book “Quran”. He does not confuses its messages. He does his research himself, resorts to conclusions he draws himself. 
p: Sect word is blatant wrong. It only creates divide. We muslims starting from Indonesia to Mauritania, all are exactly 1 thing, “Muslim”. The division among us has consumed all our energies.
h3: How do we get rid of this divide?
p: We are suffering each day. Each day is harder than the previous one. Islamic states are in a state of emergency. We need to get out of this problem sooner. How exactly can we contribute in liberating our states out of poverty?
p: We have to focus on ourselves. We have to start with “I”.
p.b-left: 1. I have to avoid listening to people who put themselves in certain boxes of religion.
p.b-left: 2. I have to read Quran more often. Quran is claimed to be an open book. Its followers do not go astray.
p.b-left: 3. I have to act upon what I read. If I get afraid I will go astray, I have to re-assure myself, it is an open book. 
p: I have to act muslim. I face an indian muslim. I embrace him. Because he believes in my Allah and my Prophet (PBUH). If I hate him because he is from India. I am still sitting in a box. There is a dire need I get out of it.
p: I have to live my life in its entirety. I can not limit myself under any circumstances; “a teacher will come to open my brain and fix me.” No teacher will come, I have to fix myself. I can not rely on any body else but my own conscience.
p: We all are miracles in ourselves. I am not better than you. You are not better than anybody else. We all are super perfect people with very logically working brain. We use this brain, we get out of trouble today. We do not use this brain, we wait for a teacher to come and fix us,  we are only delaying our exit from the box. 
h3: Should I argue with people who are in this box?
p: Do not argue. Act. Be quite. More we argue, more we lose our ground. More we act honest, more we gain the ground. Box people need teachers, imams, ancestors to move on a track. Free people need Quran to move towards the direction. 
p.note: This seems controversial but the ayat which is coming ahead, is truly a marvel. 
h3: What is the risk involved if I live in a box?
p: The Quran tells crisp and clear to live our life out of box.
quote: It is He who has sent down to you, [O Muhammad], the Book; <span class=lavendar>in it are verses [that are] precise</span> - they are the foundation of the Book - and others unspecific. <span class=lavendar>As for those in whose hearts is deviation [from truth], </span>they will follow that of it which is unspecific, seeking discord and seeking an interpretation [suitable to them].<span class=lavendar> And no one knows its [true] interpretation except Allah </span>. But those firm in knowledge say, "We believe in it. All [of it] is from our Lord." And no one will be reminded except those of understanding. <br> - Surah Al Imran verse 7
p: Quran is precise. It is really understandable for those who are firm in knowledge. People who are not so firm, focuses on unspecific, seeking discord and seek an interpretation that is suitable to them. 
p: Basic meaning of Quran is to submit to Allah (swt). It negates fitna all together. Shiaism, Sunnism, Wahabism and all sects with its sub sects are fitna if they divide muslims. People of Allah (swt) do not really need a logic to love Him. They love anybody who love Allah (swt). If on one side there is Allah (swt) on the day of judgment and on other side there is a shia brother; who we we focus on? We will focus on Allah (swt). Some shia brother, wrong practices of Islam, indian muslim will not really matter while we will be facing our Allah (swt). We will only be accountable for our relationship with Quran. 
p: Allah (swt) wants us to unite with Quran as foundation. Just holding tight to it will keep us on track. If we want leniency in some verses and we end up asking for an easier interpretation we are really going into fitna.  
h3: For example.. 
p: Hazrat Sulman Farsi (ra) was raised as a Zoroastrian, then attracted to Christianity, and then converted to Islam after meeting Muhammad in the city of Yathrib, which later became Medina. He was looking for people who were on truth before declaration of Muhammad (saw) as prophet. The basic meaning of truth was clear to him; people who act honest. His struggle to find truth lead him to Prophet (PBUH). The moment he found truth he focused on it and embraced it. 
p: After embracing Islam he was appointed as the governor of Al-Mada'in in Iraq. According to popular tradition, Muhammad considered Salman as part of his household (Ahl al-Bayt). He was a renowned follower of Ali ibn Abi Talib after the death of Muhammad. He is credited with the suggestion of digging a trench around Medina when it was attacked by Mecca in the Battle of the Trench. 
p.note: <a href=https://youtu.be/TpO89LMFcno target=_blank>Click here to listen to his story on Youtube</a>
h3: Conclusion
p: Islam is a clearly defined religion. Anything that feels wrong, really is wrong. The moment we seek leniency in its orders, we go astray. We do not know what Quran means until we act on it. 
p: Ibrahim (as)’s fire was well lit fire. People did real hard work making the biggest fire for Ibrahim (as). He was thrown into it. It was not until the catapult released him, the fire turned all green, cozy and comfy. He jumped with a blind faith into fire. Let Allah (swt) do the magic for us. Let us jump first and ask later. Let us love Him. Our fire will turn into comfort after we jump. Before jumping some molana or mufti won’t convert Allah’s fire for us to cozy, comfy and green garden. Our faith is all we are going to take with us on our death. Go blind when reading Quran. Act on it blindly. No questions asked, straight to the point, act on it for all big things only happen after we jump. 
p.note: Guys, if you like this project. Please follow this project's page on twitter. <a href="https://twitter.com/zakatlists">Click here to go to the twitter page</a>.
"@

# Row 29: new article "article 28 is live"
$ws.Range("A29").Value = 28
$ws.Range("B29").Value = 43858
$ws.Range("B29").NumberFormat = $ws.Range("B28").NumberFormat
$ws.Range("C29").Value = "Surah Baqarah, 283 - 286; Surah Al Imran, 1 - 8"
$ws.Range("F29").Value = "Understanding Quran, Trust in Allah, Avoiding arguments, Sectarian divide, Unity in Islam"
$ws.Range("D29").Value = $articleText
$ws.Range("E29").Value = "Qasim Ali"

$ws.Rows.Item(29).RowHeight = 409.6

$ws.Range("D29").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 29
$excel.ActiveWindow.ScrollColumn = 1
